{"js": "// 1. Add a hanging-indent paragraph format (left=1800 twips=90pt,\n//    hanging=360 twips=18pt -> firstLineIndent = -18pt) to the first\n//    (and only) paragraph of the main document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.leftIndent = 90; // 1800 twips / 20 = 90 points\nfirstParagraph.firstLineIndent = -18; // -(360 twips / 20) = -18 points (hanging indent)\nawait context.sync();\n\n// 2. Remove the \"PAGE \\* MERGEFORMAT\" page-number fields from all three\n//    footers (even, primary/default, first page) of the only section,\n//    leaving the ptab runs untouched.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst footerTypes = [\n  Word.HeaderFooterType.evenPages,\n  Word.HeaderFooterType.primary,\n  Word.HeaderFooterType.firstPage,\n];\n\nconst footerFieldsList = [];\nfor (const type of footerTypes) {\n  const footer = section.getFooter(type);\n  const fields = footer.fields;\n  fields.load(\"items\");\n  footerFieldsList.push(fields);\n}\nawait context.sync();\n\nfor (const fields of footerFieldsList) {\n  for (const field of fields.items) {\n    field.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# 1. Add a hanging-indent paragraph format (left=1800 twips=90pt,\n#    hanging=360 twips=18pt -> FirstLineIndent = -18pt) to the first\n#    (and only) paragraph of the main document body.\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs.Item(1)\n$firstParagraph.Range.ParagraphFormat.LeftIndent = 90        # 1800 twips / 20 = 90 points\n$firstParagraph.Range.ParagraphFormat.FirstLineIndent = -18  # -(360 twips / 20) = -18 points (hanging indent)\n\n# 2. Remove the \"PAGE \\* MERGEFORMAT\" page-number fields from all three\n#    footers (even, primary/default, first page) of the only section,\n#    leaving the ptab runs untouched.\n$section = $d.Sections.Item(1)\n\n$evenFooter = $section.Footers.Item([Microsoft.Office.Interop.Word.WdHeaderFooterIndex]::wdHeaderFooterEvenPages)\n$primaryFooter = $section.Footers.Item([Microsoft.Office.Interop.Word.WdHeaderFooterIndex]::wdHeaderFooterPrimary)\n$firstFooter = $section.Footers.Item([Microsoft.Office.Interop.Word.WdHeaderFooterIndex]::wdHeaderFooterFirstPage)\n\nforeach ($footer in @($evenFooter, $primaryFooter, $firstFooter)) {\n    $fields = $footer.Range.Fields\n    for ($i = $fields.Count; $i -ge 1; $i--) {\n        $fields.Item($i).Delete()\n    }\n}\n"}
